$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.069.42'
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("D3").Value = '2.256.48'
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.497'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.42%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.89'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0784'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.57%  '
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("D15").Value = '2.607.40'
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.60%  '
$ws.Range("D17").Value = '2.257.88'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.774'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.24%  '
$ws.Range("D19").Value = '42.070.77'
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("D20").Value = '0.0₃0886'
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.55%  '
$ws.Range("E22").Value = '  -4.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '232.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("E25").Value = '  -4.50%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.16%  '
$ws.Range("E30").Value = '  -12.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.81%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("E34").Value = '  -3.98%  '
$ws.Range("E35").Value = '  -4.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0691'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.15%  '
$ws.Range("E37").Value = '  -5.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0982'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.04%  '
$ws.Range("E41").Value = '  -3.28%  '
$ws.Range("E42").Value = '  -8.80%  '
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").Value = '1.933.15'
$ws.Range("E44").Value = '  -4.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0278'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.28%  '
$ws.Range("E48").Value = '  -5.57%  '
$ws.Range("E49").Value = '  -2.94%  '
$ws.Range("D50").Value = '2.482.56'
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.40%  '
